$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: Expand author initials to full first names in the opening citation.
# ---------------------------------------------------------------------------
$find = $d.Content.Find
$find.Execute(
    "Junker, J. R., W. F. Cross, J. M. Hood, J. P. Benstead, A. D. Huryn, D. Nelson, J. S. Ólafsson, and G. M. Gíslason,",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "James R. Junker, Wyatt F. Cross, James M. Hood, Jonathan P. Benstead, Alexander D. Huryn, Daniel Nelson, Jón S. Ólafsson, and Gísli M. Gíslason,",
    2) | Out-Null

Write-Host "Edit 1 done"

# ---------------------------------------------------------------------------
# Edit 2: Figure S1 image caption - the sentence is re-flowed into two runs
# split at a different point (text content is identical).
# ---------------------------------------------------------------------------
$rng = $d.Content.Duplicate
$rng.Find.Execute("Figure S1. Daily mean temperature") | Out-Null
$capPara = $rng.Paragraphs(1).Range

$xml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ImageCaption"/></w:pPr><w:r><w:t>Figure S1. Daily mean temperature (°C) across study streams for each day of year (doy) over the course of th</w:t></w:r><w:r><w:t>e study. Legend represents the annual mean temperate within each stream. This figure was modified from Junker 2019 with permissiom.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
$capPara.InsertXML($xml)

Write-Host "Edit 2 done"

# ---------------------------------------------------------------------------
# Edit 3: Figure S2 caption - re-flow two runs at a different split point and
# drop the trailing manual page-break run (the whole paragraph, including the
# embedded equation run, is rewritten verbatim).
# ---------------------------------------------------------------------------
$rng2 = $d.Content.Duplicate
$rng2.Find.Execute("Appendix S2: Figure S2. Species-level distributions") | Out-Null
$fig2Para = $rng2.Paragraphs(1).Range

$xml2 = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:body><w:p><w:pPr><w:pStyle w:val="ImageCaption"/></w:pPr><w:r><w:t xml:space="preserve">Appendix S2: Figure S2. Species-level distributions of population variables (median </w:t></w:r><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>±</m:t></m:r></m:oMath><w:r><w:t xml:space="preserve"> median abso</w:t></w:r><w:r><w:t>lute deviations) observed across all streams in the study. a) total annual organic matter flux, b) standing population biomass, c) population abundance, and d) mean individual body size. All variables have been log</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="subscript"/></w:rPr><w:t>10</w:t></w:r><w:r><w:t>-transformed.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
$fig2Para.InsertXML($xml2)

Write-Host "Edit 3 done"

# ---------------------------------------------------------------------------
# Edit 4: "Appendix S2: Table S1" heading - merge the two runs back into one,
# and relocate the (Word-managed) "_GoBack" bookmark here, collapsed, right
# before the text run - this is where the cursor was left before saving.
# ---------------------------------------------------------------------------
$rng3 = $d.Content.Duplicate
$rng3.Find.Execute("Appendix S2: Table S1") | Out-Null
$tableS1Para = $rng3.Paragraphs(1).Range

$xml3 = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Appendix S2: Table S1</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
$tableS1Para.InsertXML($xml3)

$rng3b = $d.Content.Duplicate
$rng3b.Find.Execute("Appendix S2: Table S1") | Out-Null
$goBackSpot = $d.Range($rng3b.Start, $rng3b.Start)
$d.Bookmarks.Add("_GoBack", $goBackSpot) | Out-Null

Write-Host "Edit 4 done"

# ---------------------------------------------------------------------------
# Edit 5: Table S1 caption - re-flow the two runs at a different split point
# (text content identical).
# ---------------------------------------------------------------------------
$rng4 = $d.Content.Duplicate
$rng4.Find.Execute("Appendix S2:Table S1. Evenness") | Out-Null
$tblCaptionPara = $rng4.Paragraphs(1).Range

$xml4 = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="TableCaption"/></w:pPr><w:r><w:t>Ap</w:t></w:r><w:r><w:t>pendix S2:Table S1. Evenness of organic matter fluxes among consumers within a stream community measured by the Gini index, both raw ('non-normalized') and 'normalized' for consumer richness</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
$tblCaptionPara.InsertXML($xml4)

Write-Host "Edit 5 done"

# ---------------------------------------------------------------------------
# Edit 6: Table cell "0.22 ( 0.18 - 0.27 )" - split the trailing " )" off
# into its own run.
# ---------------------------------------------------------------------------
$rng5 = $d.Content.Duplicate
$rng5.Find.Execute("0.22 ( 0.18 - 0.27 )") | Out-Null
$cellPara = $rng5.Paragraphs(1).Range

$xml5 = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="right"/></w:pPr><w:r><w:t>0.22 ( 0.18 - 0.27</w:t></w:r><w:r><w:t xml:space="preserve"> )</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
$cellPara.InsertXML($xml5)

Write-Host "Edit 6 done"

# ---------------------------------------------------------------------------
# Edit 7: Bibliography entry - drop the direct-formatting overrides on the
# paragraph (spacing/indent/font) and on its two runs, reverting them to the
# plain "Bibliography" style; also drop the now-orphaned "_GoBack" bookmark
# that used to live here (it was relocated in Edit 4).
# ---------------------------------------------------------------------------
$rng6 = $d.Content.Duplicate
$rng6.Find.Execute("Junker, J. R. 2019, November") | Out-Null
$biblioPara = $rng6.Paragraphs(1).Range

# This is the very last paragraph in the document body, so InsertXML on its
# own range leaves a trailing stub paragraph behind (no paragraph mark to
# reclaim past the end of the story). Work around it by also swallowing the
# preceding paragraph mark and re-emitting the "References" heading verbatim,
# then deleting the leftover empty stub afterwards.
$prevParaEnd = $biblioPara.Start
$replaceRange = $d.Range($prevParaEnd - 1, $biblioPara.End)

$xml6 = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Bibliography"/></w:pPr><w:bookmarkStart w:id="100" w:name="ref-junker2019"/><w:bookmarkStart w:id="101" w:name="refs"/><w:r><w:t>Junker, J. R. 2019, November. The</w:t></w:r><w:r><w:t xml:space="preserve"> effects of temperature on stream ecosystem structure, secondary production, and food web dynamics. Doctoral, Montana State University, Bozeman, MT.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
$replaceRange.InsertXML($xml6)

$lastPara = $d.Paragraphs($d.Paragraphs.Count)
if ($lastPara.Range.Text -eq "") {
    $prevPara = $d.Paragraphs($d.Paragraphs.Count - 1)
    $cleanup = $d.Range($prevPara.Range.End - 1, $lastPara.Range.End)
    $cleanup.Delete()
}

Write-Host "Edit 7 done"
